$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain text (not parseable as a pure number)
# are set directly via .Value, which Excel stores as text.
$ws.Range("D2").Value = "26.150.18"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").Value = "1.655.46"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("E6").Value = "  -2.39%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("E10").Value = "  -1.29%  "
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.662.25"
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").Value = "1.882.94"
$ws.Range("E14").Value = "  -1.62%  "
$ws.Range("E15").Value = "  -2.90%  "
$ws.Range("D16").Value = "0.0₅8244"
$ws.Range("E16").Value = "  -2.10%  "
$ws.Range("E17").Value = "  -1.89%  "
$ws.Range("D18").Value = "26.206.52"
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("E20").Value = "  -3.13%  "
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("E22").Value = "  -2.51%  "
$ws.Range("E23").Value = "  -4.58%  "
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("E25").Value = "  -4.66%  "
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("E27").Value = "  -3.32%  "
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("E32").Value = "  -1.46%  "
$ws.Range("E33").Value = "  -3.43%  "
$ws.Range("E34").Value = "  -2.80%  "
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("E38").Value = "  +3.73%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").Value = "1.036.43"
$ws.Range("E43").Value = "  -1.53%  "
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").Value = "1.797.65"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈109"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("E51").Value = "  +3.91%  "

# Cells whose new values look like plain numbers (e.g. "218.59") must be
# forced to remain text, matching the original text-cell type, by applying
# a Text number format before assigning the value.
$numericTextCells = @{
    "D5" = "218.59"
    "D6" = "0.5208"
    "D8" = "0.2670"
    "D9" = "0.06330"
    "D10" = "21.08"
    "D11" = "0.07754"
    "D13" = "4.437"
    "D15" = "0.5463"
    "D20" = "4.669"
    "D21" = "192.92"
    "D22" = "10.16"
    "D23" = "6.103"
    "D25" = "137.42"
    "D26" = "0.1237"
    "D27" = "7.239"
    "D28" = "16.12"
    "D29" = "1.410"
    "D30" = "0.06041"
    "D31" = "1.286"
    "D32" = "3.555"
    "D33" = "3.343"
    "D34" = "1.653"
    "D35" = "0.9821"
    "D36" = "2.409"
    "D38" = "0.5930"
    "D39" = "0.01595"
    "D40" = "5.965"
    "D41" = "0.8659"
    "D42" = "1.003"
    "D44" = "99.89"
    "D47" = "57.21"
    "D48" = "1.007"
    "D49" = "8.100"
    "D50" = "0.05178"
    "D51" = "1.471"
}
foreach ($ref in $numericTextCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextCells[$ref]
}
